# Fixed update to excel issue
# - Rename "Requested quantity" headers to metric-specific names
# - Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

# --- 1. Rename headers on the existing sheets -------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the exact same cell formatting used on the other sheets:
#   header row style (bold, centered, bordered) and date-column style (datetime number format)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A18").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Forecast data rows (row 2 .. row 18)
$data = @(
    @(44948.99999999999, 33, 7.697793126795123, 60.98789767679598),
    @(44955.99999999999, 32, 5.586601351227396, 58.30975446082111),
    @(44962.99999999999, 30, 4.906990838441518, 57.29340622818079),
    @(44976.99999999999, 26, -1.609844408791466, 52.97224007781583),
    @(44983.99999999999, 24, -1.577531808025765, 52.95985862577633),
    @(44990.99999999999, 22, -5.022593463005325, 49.43066994692224),
    @(44997.99999999999, 21, -6.22212182129167, 48.3152052593071),
    @(45004.99999999999, 19, -7.330646659153258, 45.04438145740706),
    @(45011.99999999999, 17, -10.94516472273024, 44.63269900383381),
    @(45018.99999999999, 15, -11.19338691460407, 41.00823727168792),
    @(45025.99999999999, 13, -12.88565583643481, 39.49294991095252),
    @(45032.99999999999, 11, -16.24945829322794, 38.96423435473208),
    @(45039.99999999999, 10, -17.42130884134346, 37.37362121812424),
    @(45046.99999999999, 8, -19.58229997508099, 35.27739668342043),
    @(45053.99999999999, 6, -20.19448886239783, 34.69063146932245),
    @(45060.99999999999, 4, -24.4542985045691, 32.52585943961298),
    @(45067.99999999999, 2, -25.90301459353765, 29.70582709506111)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$wsForecast.Range("A1").Select()
